$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Names (A38:A40)
$ws.Range("A38").Value = "MAHESWARI P"
$ws.Range("A39").Value = "SASI KALA.B"
$ws.Range("A40").Value = "SWETHA K R"

# Designation (B38:B40) - same text each time
$ws.Range("B38").Value = "Lab Instructor"
$ws.Range("B39").Value = "Lab Instructor"
$ws.Range("B40").Value = "Lab Instructor"

# Photo URLs (C38:C40) written in ascending numeric order 20,21,22
$ws.Range("C40").Value = "/static/images/profile_photos/005/VEC-005-05-20.webp"
$ws.Range("C39").Value = "/static/images/profile_photos/005/VEC-005-05-21.webp"
$ws.Range("C38").Value = "/static/images/profile_photos/005/VEC-005-05-22.webp"

# Unique IDs (J38:J40) in row order 22,21,20
$ws.Range("J38").Value = "VEC-005-05-22"
$ws.Range("J39").Value = "VEC-005-05-21"
$ws.Range("J40").Value = "VEC-005-05-20"

# Column A on the new rows gets a thin box border (matches the style used
# when Excel's "All Borders" is applied from the Home ribbon).
$ws.Range("A38:A40").Borders.LineStyle = 1

# Restore the view/selection the author left the sheet in.
$ws.Range("I38").Select()

Write-Host "done"
